$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8. This shifts all existing data rows
# (8..99) down by one (to 9..100), matching the diff which shows every
# existing record moving one row down and a new record appearing at the
# very bottom (row 100, formerly row 99's data).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44545
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región de La Araucanía"
$ws.Range("P8").Value = 800
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"

# Preserve the date number format used by the rest of column D.
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
